# Membuat fitur alokasi mitra
# Rename the "idsobat" header (A1) to "sobat_id" and move it to column A1,
# keeping the rest of the header row the same. Also move the active
# selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "sobat_id"

$ws.Range("A2").Select()
